# Generate Report for Handback
#
# The localization job for "139671fe-2709-4645-a438-38e543dc5459.md" has
# been handed back (translations are in sync with en-US). This script
# updates the Overview sheet and the per-language detail sheets (zh-cn,
# de-de) to reflect the handback: new status text, a populated "Latest
# Target File" / "Latest Handback File" pair (with hyperlinks mirroring
# the existing Source/Handoff links), a handback timestamp, and the
# "Include" handoff reason.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status columns for the file that
# was handed back.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Helper: update one language detail sheet (zh-cn / de-de) for row 2
# (the 139671fe... file), given the new "Latest Handback DateTime".
# ---------------------------------------------------------------------
function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Find the existing hyperlink addresses for A2 (Source File Name) and
    # D2 (Latest Handoff File) so the new F2/G2 hyperlinks can mirror them.
    $sourceAddr = $null
    $handoffAddr = $null
    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range
        if ($r.Row -eq 2 -and $r.Column -eq 1) { $sourceAddr = $hl.Address }
        if ($r.Row -eq 2 -and $r.Column -eq 4) { $handoffAddr = $hl.Address }
    }

    # Status -> Handed back
    $ws.Range("C2").Value = $statusHandedBack

    # Latest Target File (F2) - mirrors the Source File Name hyperlink
    $ws.Hyperlinks.Add($ws.Range("F2"), $sourceAddr, "", "", "139671fe-2709-4645-a438-38e543dc5459.md") | Out-Null

    # Latest Handback File (G2) - mirrors the Latest Handoff File hyperlink
    $ws.Hyperlinks.Add($ws.Range("G2"), $handoffAddr, "", "", $ws.Range("D2").Value) | Out-Null

    # Latest Handback DateTime (H2)
    $ws.Range("H2").Value = $HandbackDateTime

    # Handoff Reason (I2)
    $ws.Range("I2").Value = "Include"
}

Update-HandbackSheet "zh-cn" "2016-03-18 20:32:34"
Update-HandbackSheet "de-de" "2016-03-18 20:32:39"
